$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Form Responses 1")

$data = @(
    @("12/31/2020 2:54:16", "b20108@students.iitmandi.ac.in", "Kanchan Padvi", "https://www.urionlinejudge.com.br/judge/en/profile/509219", "B20108", "ComputerScience"),
    @("12/30/2020 21:24:43", "b20133@students.iitmandi.ac.in", "Shailesh Rathod", "https://www.urionlinejudge.com.br/judge/en/profile/509167", "B20133", "ComputerScience"),
    @("12/31/2020 1:02:28", "b20124@students.iitmandi.ac.in", "B20124 Rajeev Kumar", "https://www.urionlinejudge.com.br/judge/en/profile/509266", "B20124", "ComputerScience"),
    @("12/31/2020 11:46:59", "b20162@students.iitmandi.ac.in", "Pavitra Jain", "https://www.urionlinejudge.com.br/judge/en/profile/509047", "B20162", "DataScience"),
    @("12/31/2020 13:50:59", "b20097@students.iitmandi.ac.in", "DIYA ASHISH ", "https://www.urionlinejudge.com.br/judge/en/profile/508228", "B20097", "ComputerScience")
)

$startRow = 100
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $dCell = $ws.Cells.Item($r, 4)
    $dCell.Value = $row[3]
    $null = $ws.Hyperlinks.Add($dCell, $row[3], "", "", $row[3])
    $dCell.Font.Name = "Arial"
    $dCell.Font.Size = 10
    $dCell.Font.Underline = 0
    $dCell.Font.Color = 0
    $dCell.Font.Bold = $false
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}

$null = $ws.Range("E101").Select()
